$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data (avoids Excel
# auto-converting numeric-looking strings like "305.83" into real numbers).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Cell value updates ---
$ws.Range('D2').Value = '44.053.35'
$ws.Range('D3').Value = '2.242.27'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '305.83'
$ws.Range('E5').Value = '  -5.19%  '
$ws.Range('D6').Value = '96.60'
$ws.Range('E6').Value = '  -5.13%  '
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  -5.50%  '
$ws.Range('D10').Value = '34.74'
$ws.Range('E10').Value = '  -6.39%  '
$ws.Range('E11').Value = '  -2.95%  '
$ws.Range('D12').Value = '7.17'
$ws.Range('E12').Value = '  -7.18%  '
$ws.Range('D14').Value = '2.584.87'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = '2.240.30'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '0.822'
$ws.Range('E16').Value = '  -4.19%  '
$ws.Range('E17').Value = '  -4.19%  '
$ws.Range('D18').Value = '43.884.99'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').Value = '0.0₃0966'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').Value = '12.36'
$ws.Range('E20').Value = '  -9.19%  '
$ws.Range('E21').Value = '  -5.07%  '
$ws.Range('D22').Value = '64.78'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = '238.91'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = '2.93'
$ws.Range('E24').Value = '  -7.79%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  -9.81%  '
$ws.Range('D27').Value = '9.93'
$ws.Range('E27').Value = '  -2.48%  '
$ws.Range('E28').Value = '  -3.13%  '
$ws.Range('D29').Value = '36.36'
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('D30').Value = '6.03'
$ws.Range('E30').Value = '  -4.61%  '
$ws.Range('D31').Value = '19.98'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('D32').Value = '153.63'
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('D33').Value = '3.36'
$ws.Range('E33').Value = '  +9.41%  '
$ws.Range('D34').Value = '0.0806'
$ws.Range('E34').Value = '  -5.58%  '
$ws.Range('D35').Value = '2.66'
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('D36').Value = '0.119'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  -7.07%  '
$ws.Range('D38').Value = '1.76'
$ws.Range('E38').Value = '  -9.24%  '
$ws.Range('D39').Value = '14.83'
$ws.Range('E39').Value = '  -7.05%  '
$ws.Range('D40').Value = '3.79'
$ws.Range('E40').Value = '  -10.82%  '
$ws.Range('D41').Value = '0.0302'
$ws.Range('E41').Value = '  -5.10%  '
$ws.Range('E42').Value = '  -12.50%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').Value = '1.748.63'
$ws.Range('E44').Value = '  -3.63%  '
$ws.Range('D45').Value = '85.86'
$ws.Range('E45').Value = '  +3.96%  '
$ws.Range('D46').Value = '16.01'
$ws.Range('E46').Value = '  +13.12%  '
$ws.Range('D47').Value = '5.10'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('E48').Value = '  -5.77%  '
$ws.Range('D49').Value = '100.37'
$ws.Range('E49').Value = '  -3.15%  '
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = '68.69'
$ws.Range('E51').Value = '  -9.79%  '

# Restore default style on column D so no stray per-cell formatting remains.
$priceRange.Style = "Normal"
